# Apply text replacements to update the document date and the division expressions
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-13 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-14 Thursday", 2) | Out-Null
$d.Content.Find.Execute("768÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "734÷2=", 2) | Out-Null
$d.Content.Find.Execute("570÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "386÷4=", 2) | Out-Null
$d.Content.Find.Execute("929÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "829÷7=", 2) | Out-Null
$d.Content.Find.Execute("210÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "937÷9=", 2) | Out-Null
$d.Content.Find.Execute("409÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "913÷7=", 2) | Out-Null
$d.Content.Find.Execute("135÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "898÷8=", 2) | Out-Null
$d.Content.Find.Execute("242÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "590÷3=", 2) | Out-Null
$d.Content.Find.Execute("213÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "337÷8=", 2) | Out-Null
$d.Content.Find.Execute("291÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "155÷8=", 2) | Out-Null
$d.Content.Find.Execute("412÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "573÷3=", 2) | Out-Null
$d.Content.Find.Execute("960÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "569÷3=", 2) | Out-Null
$d.Content.Find.Execute("516÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "436÷7=", 2) | Out-Null
$d.Content.Find.Execute("718÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "256÷3=", 2) | Out-Null
$d.Content.Find.Execute("186÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "244÷6=", 2) | Out-Null
$d.Content.Find.Execute("607÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "441÷8=", 2) | Out-Null
$d.Content.Find.Execute("355÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "693÷7=", 2) | Out-Null
$d.Content.Find.Execute("840÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "632÷4=", 2) | Out-Null
$d.Content.Find.Execute("405÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "125÷3=", 2) | Out-Null
$d.Content.Find.Execute("379÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "529÷9=", 2) | Out-Null
$d.Content.Find.Execute("963÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "857÷2=", 2) | Out-Null
$d.Content.Find.Execute("890÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "142÷5=", 2) | Out-Null
$d.Content.Find.Execute("316÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "848÷9=", 2) | Out-Null
$d.Content.Find.Execute("684÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "449÷3=", 2) | Out-Null
$d.Content.Find.Execute("176÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "475÷8=", 2) | Out-Null
$d.Content.Find.Execute("341÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "983÷6=", 2) | Out-Null
